$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume snapshot (scraped data update).
# Column D holds formatted price strings (e.g. "216.90", "1.01") that must
# stay literal text rather than being auto-coerced to numbers (which would
# drop trailing zeros / reformat the value) -- force text via a leading
# apostrophe on .Formula, matching how the source data is authored.

$ws.Range('D2').Formula = "'26.954.64"
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Formula = "'1.563.62"
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('D4').Formula = "'1.01"
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Formula = "'207.74"
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Formula = "'22.09"
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('D10').Formula = "'0.0601"
$ws.Range('E10').Value = '  +2.13%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Formula = "'1.785.00"
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Formula = "'1.563.05"
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Formula = "'26.953.20"
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').Formula = "'216.90"
$ws.Range('E19').Value = '  -1.05%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').Formula = "'152.71"
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').Formula = "'15.07"
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').Value = '  +0.99%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').Formula = "'1.12"
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').Formula = "'3.12"
$ws.Range('E33').Value = '  +1.16%  '
$ws.Range('D34').Formula = "'1.422.50"
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('E35').Value = '  +2.87%  '
$ws.Range('E36').Value = '  +9.76%  '
$ws.Range('D37').Formula = "'2.34"
$ws.Range('E37').Value = '  +1.69%  '
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('D39').Formula = "'0.535"
$ws.Range('E39').Value = '  +1.67%  '
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('E41').Value = '  -0.88%  '
$ws.Range('D42').Formula = "'1.01"
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('D45').Formula = "'64.86"
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('D47').Formula = "'1.698.95"
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Formula = "'87.44"
$ws.Range('E48').Value = '  +0.63%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Formula = "'0.0₆0102"
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Formula = "'0.0519"
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Formula = "'0.0959"
$ws.Range('E51').Value = '  -0.74%  '
